$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "ClassDiagram" to "class Patient"
$ws.Name = "class Patient"

# Add the new attribute row: "int" | "respiration_rate"
$ws.Range("A2").Value = "int"
$ws.Range("B2").Value = "respiration_rate"

# Center the new attribute-type column (column A). Using the full column
# range (rather than just Columns.Item) ensures the newly written cell
# A2 correctly picks up the centered, "Aptos Display" bodied style that
# already backs column A/B, matching the style used for column B's font.
$ws.Range("A1:A1048576").HorizontalAlignment = -4108

# Give column A (the type column) its own, narrower width, distinct from
# column B's width, matching the layout used for the new table row.
$ws.Columns.Item(1).ColumnWidth = 23.8333333333

# Move the active selection to A3, below the newly added row.
$ws.Range("A3").Select()
